$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45097
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("S2").Value = 944

# Row 3
$ws.Range("D3").Value = 44299
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 15

# Row 4
$ws.Range("D4").Value = 44294
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 800

# Row 5
$ws.Range("D5").Value = 44319
$ws.Range("M5").Value = 120

# Row 6
$ws.Range("D6").Value = 44342
$ws.Range("L6").Value = 'Especial'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 1111

# Row 7
$ws.Range("D7").Value = 44355
$ws.Range("L7").Value = 'Especial'
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44340
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 230

# Row 9
$ws.Range("D9").Value = 44326
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 1111
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 45099
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("S10").Value = 944

# Row 11
$ws.Range("D11").Value = 44691

# Row 12
$ws.Range("D12").Value = 45083
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 17000
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 944

# Row 13
$ws.Range("D13").Value = 45085
$ws.Range("M13").Value = 280
$ws.Range("N13").Value = 17000
$ws.Range("P13").Value = 17357
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 964

# Row 14
$ws.Range("D14").Value = 44316
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("S14").Value = 1111

# Row 15
$ws.Range("D15").Value = 45062
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("Q15").Value = '$/caja 15 kilos granel'
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 15

# Row 16
$ws.Range("D16").Value = 44328
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1111

# Row 17
$ws.Range("D17").Value = 44692
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 150

# Row 18
$ws.Range("D18").Value = 44714
$ws.Range("M18").Value = 100

# Row 19
$ws.Range("D19").Value = 44348
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 200

# Row 20
$ws.Range("D20").Value = 44358
$ws.Range("N20").Value = 18000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 18000
$ws.Range("S20").Value = 1000

# Row 21
$ws.Range("D21").Value = 44358
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 17000
$ws.Range("P21").Value = 17000
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 944
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44354
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 18000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 18000
$ws.Range("Q22").Value = '$/caja 18 kilos granel'
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44680
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 15000
$ws.Range("Q23").Value = '$/caja 15 kilos granel'
$ws.Range("R23").Value = 'Provincia de Limarí'
$ws.Range("S23").Value = 1000
$ws.Range("T23").Value = 15

# Row 24
$ws.Range("D24").Value = 44291
$ws.Range("M24").Value = 150
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("Q24").Value = '$/caja 15 kilos granel'
$ws.Range("R24").Value = 'Región Metropolitana'
$ws.Range("S24").Value = 800
$ws.Range("T24").Value = 15

# Row 25
$ws.Range("D25").Value = 45054
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("S25").Value = 778
